$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/1h (E) values for rows with changed figures.
# Column D assignments use a leading apostrophe (via .Formula) for values that
# would otherwise be auto-parsed by Excel as numbers (e.g. 406.39, 4.00), which
# keeps them as text and avoids floating point rounding artifacts / lost trailing zeros.

$ws.Range("D2").Value = "61.805.07"
$ws.Range("E2").Value = "  +2.17%  "

$ws.Range("D3").Value = "3.403.17"
$ws.Range("E3").Value = "  +3.95%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Formula = "'406.39"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").Formula = "'130.74"
$ws.Range("E6").Value = "  +18.25%  "

$ws.Range("D7").Formula = "'0.610"
$ws.Range("E7").Value = "  +8.47%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Formula = "'0.676"
$ws.Range("E9").Value = "  +10.23%  "

$ws.Range("D10").Formula = "'0.127"
$ws.Range("E10").Value = "  +13.07%  "

$ws.Range("D11").Formula = "'42.40"
$ws.Range("E11").Value = "  +10.70%  "

$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").Value = "3.956.65"
$ws.Range("E13").Value = "  +3.77%  "

$ws.Range("D14").Formula = "'8.56"
$ws.Range("E14").Value = "  +6.30%  "

$ws.Range("E15").Value = "  +4.92%  "

$ws.Range("D16").Value = "3.395.30"
$ws.Range("E16").Value = "  +1.77%  "

$ws.Range("D17").Formula = "'11.51"
$ws.Range("E17").Value = "  +10.22%  "

$ws.Range("D18").Value = "61.580.33"
$ws.Range("E18").Value = "  +1.78%  "

$ws.Range("E19").Value = "  +5.21%  "

$ws.Range("E20").Value = "  +20.03%  "

$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("D22").Formula = "'82.64"
$ws.Range("E22").Value = "  +13.60%  "

$ws.Range("D23").Formula = "'13.15"
$ws.Range("E23").Value = "  +6.15%  "

$ws.Range("D24").Formula = "'308.26"
$ws.Range("E24").Value = "  +4.51%  "

$ws.Range("D25").Formula = "'3.18"
$ws.Range("E25").Value = "  +3.66%  "

$ws.Range("D26").Formula = "'8.59"
$ws.Range("E26").Value = "  +16.20%  "

$ws.Range("E27").Value = "  +3.25%  "

$ws.Range("D28").Formula = "'4.68"
$ws.Range("E28").Value = "  +9.77%  "

$ws.Range("D29").Formula = "'7.47"
$ws.Range("E29").Value = "  +2.51%  "

$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Formula = "'11.85"
$ws.Range("E31").Value = "  +7.23%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Formula = "'0.116"
$ws.Range("E32").Value = "  +4.15%  "

$ws.Range("D33").Formula = "'2.61"
$ws.Range("E33").Value = "  +7.97%  "

$ws.Range("D34").Formula = "'42.48"
$ws.Range("E34").Value = "  +10.25%  "

$ws.Range("D35").Formula = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  +3.05%  "

$ws.Range("D37").Formula = "'52.38"
$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("E38").Value = "  -0.61%  "

$ws.Range("D39").Formula = "'3.44"
$ws.Range("E39").Value = "  +5.57%  "

$ws.Range("D40").Formula = "'2.99"
$ws.Range("E40").Value = "  -2.96%  "

$ws.Range("E41").Value = "  +9.78%  "

$ws.Range("E42").Value = "  +5.95%  "

$ws.Range("D43").Formula = "'136.81"
$ws.Range("E43").Value = "  +1.95%  "

$ws.Range("D44").Formula = "'4.00"
$ws.Range("E44").Value = "  +8.11%  "

$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Formula = "'17.12"
$ws.Range("E45").Value = "  +6.52%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Formula = "'0.286"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("E47").Value = "  +1.61%  "

$ws.Range("D48").Formula = "'21.83"
$ws.Range("E48").Value = "  +5.33%  "

$ws.Range("D49").Value = "2.148.97"
$ws.Range("E49").Value = "  +2.37%  "

$ws.Range("D50").Value = "3.739.64"
$ws.Range("E50").Value = "  +3.12%  "

$ws.Range("D51").Formula = "'2.35"
$ws.Range("E51").Value = "  +0.38%  "

